$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: row, column C (nombre_aides), column E (montant_total)
$updates = @(
    @{Row=3;   C=249327;  E=1036476130},
    @{Row=4;   C=103466;  E=627499207},
    @{Row=62;  C=4189;    E=9186399},
    @{Row=91;  C=151105;  E=482131030},
    @{Row=92;  C=409019;  E=1593716220},
    @{Row=93;  C=209507;  E=1307872027},
    @{Row=94;  C=94148;   E=915830612},
    @{Row=95;  C=50728;   E=930867754},
    @{Row=96;  C=17247;   E=789598588},
    @{Row=104; C=135223;  E=272132060},
    @{Row=114; C=3799;    E=9102060},
    @{Row=115; C=11692;   E=32953239},
    @{Row=116; C=4560;    E=20567318},
    @{Row=127; C=57;      E=207876},
    @{Row=131; C=75581;   E=307211165},
    @{Row=157; C=12097;   E=182911050}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
